# Requirements.xlsx update (issue #332 - update documentation due to modularization)
#
# Content change: the "General" sheet's G-6 requirement ("available via public
# repositories") drops the outdated description referencing the defunct
# bintray/jcenter repositories and replaces it with "Maven Central" (the
# project dropped JCenter in favor of Maven Central / JitPack).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

# G-6 description: "like bintray or jcenter." -> "Maven Central"
$ws.Range("C7").Value = "Maven Central"

# The cells below only carried redundant/unused direct formatting that has no
# visible effect (equivalent to the sheet's default style); drop it so the
# cells fall back to the default style, same as after the source file was
# re-saved.
$ws.Range("B3").ClearFormats()
$ws.Range("C3").ClearFormats()
$ws.Range("C6").ClearFormats()

# Update the active selection on the "General" tab to reflect where editing
# left off.
$ws.Activate()
$ws.Range("C7").Select()
